# Generate Report for Handoff
#
# File "b.md" moves from "Handed back: in sync with en-US" to
# "Ready for handoff": a new handoff package (b.<hash>.<locale>.xlf) was
# generated, and the localization-status report is regenerated to reflect
# it (plus a warning that the handback package on file is stale relative
# to the newly generated handoff).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-12 08:43:57"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-12 08:43:51"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/662b50cd2fc596b7efb057adf6abdb57d2af64e0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/1675b20c8366eba72185572b63aa67fb7aa1a5ad/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-12 08:43:57"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/662b50cd2fc596b7efb057adf6abdb57d2af64e0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/1675b20c8366eba72185572b63aa67fb7aa1a5ad/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
